# "Add files via upload" — the transactions export was refreshed:
#   * the rekey transaction (old row 6, and its "rekey"/"from"/"to" columns
#     G:I) was dropped entirely
#   * the asset-config row (old row 5) now documents a different ASA
#     (ESG6 / eresung6) with updated addr/clawback/manager values
#   * selection moved to the new last data cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the rekey transaction's dedicated columns (old G:I => rekey/from/to).
# Everything to the right (assset_type, addr, defaultFrozen, ...) shifts left.
$ws.Range("G:I").Delete()

# Drop the rekey transaction row itself (old row 6).
$ws.Range("6:6").Delete()

# Refresh the asset-config details (row 5) to describe the new asset.
$ws.Range("H5").Value = "GGY4WEN5FLISHFSXEPAELC6JWFB7R3UVTB3DWBPFUCBXZ3J26I2LDVCYZM"
$ws.Range("N5").Value = "2WEXHKWRYK6MQLNTW7GMDX72ZKAW3TVJ52SZ2EDAE23OF2Q34UUOGGT7CY"
$ws.Range("O5").Value = "2WEXHKWRYK6MQLNTW7GMDX72ZKAW3TVJ52SZ2EDAE23OF2Q34UUOGGT7CY"
$ws.Range("Q5").Value = "eresung6"
$ws.Range("P5").Value = "ESG6"

# Match the author's final selection/scroll state.
$ws.Range("P5").Select()
$excel.ActiveWindow.ScrollColumn = 9
